$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '24.582.82'
$ws.Range('E2').Value = '  -0.24%  '
$ws.Range('D3').Value = '1.685.16'
$ws.Range('E3').Value = '  -0.21%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '1.009'
$c.Style = "Normal"
$ws.Range('E4').Value = '  +0.58%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '314.26'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.49%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.43%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.3897'
$c.Style = "Normal"
$ws.Range('E7').Value = '  -1.01%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.4032'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.18%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '1.488'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -0.11%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '1.008'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.43%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '53.09'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +0.28%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.08698'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.05%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '7.586'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +5.31%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '24.64'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +5.17%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '7.918'
$c.Style = "Normal"
$ws.Range('E15').Value = '  -1.62%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '0.00001336'
$c.Style = "Normal"
$ws.Range('E16').Value = '  +1.88%  '
$ws.Range('D17').Value = '1.689.96'
$ws.Range('E17').Value = '  -0.41%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '98.27'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -1.29%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.07089'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +1.33%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '19.66'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +1.34%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '7.292'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +4.06%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '1.010'
$c.Style = "Normal"
$ws.Range('E22').Value = '  +0.80%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '14.16'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.66%  '
$ws.Range('D24').Value = '24.653.46'
$ws.Range('E24').Value = '  +0.12%  '
$ws.Range('E25').Value = '  -6.97%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '2.345'
$c.Style = "Normal"
$ws.Range('E26').Value = '  -0.84%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '22.65'
$c.Style = "Normal"
$ws.Range('E27').Value = '  +0.26%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '161.54'
$c.Style = "Normal"
$ws.Range('E28').Value = '  -0.76%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '8.442'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +11.40%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '5.236'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +0.93%  '
$c = $ws.Range('D31')
$c.NumberFormat = "@"
$c.Value = '136.41'
$c.Style = "Normal"
$ws.Range('E31').Value = '  +0.95%  '
$ws.Range('D32').Value = '1.874.52'
$ws.Range('E32').Value = '  -0.45%  '
$ws.Range('B33').Value = 'InternetComputer(DFINITY)'
$ws.Range('C33').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '7.564'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +5.16%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.08785'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +2.93%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.033'
$c.Style = "Normal"
$ws.Range('E35').Value = '  -1.62%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '1.984'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +4.34%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.02911'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +7.73%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.2714'
$c.Style = "Normal"
$ws.Range('E38').Value = '  -0.39%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '10.69'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -3.96%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.09090'
$c.Style = "Normal"
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('B41').Value = 'Aptos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '14.05'
$c.Style = "Normal"
$ws.Range('E41').Value = '  -1.56%  '
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '0.7810'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +3.08%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '1.458'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +0.09%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '16.61'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +4.38%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.7142'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +0.44%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.568'
$c.Style = "Normal"
$ws.Range('E46').Value = '  -0.21%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '4.200'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -0.20%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '1.006'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +0.47%  '
$ws.Range('E49').Value = '  +1.30%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '138.04'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -0.57%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '90.45'
$c.Style = "Normal"
$ws.Range('E51').Value = '  +1.21%  '
